$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (A: date serial, B: nuovi pos., C: somma mobile 7gg., D: somma mobile 7gg. per 100mila abitanti)
$data = @(
  @(44403, 1, 6, 53.14908317831517),
  @(44404, 0, 6, 53.14908317831517),
  @(44405, 0, 6, 53.14908317831517),
  @(44406, 2, 8, 70.86544423775356),
  @(44407, 5, 11, 97.43998582691115),
  @(44408, 0, 8, 70.86544423775356),
  @(44409, 1, 9, 79.72362476747276),
  @(44410, 0, 8, 70.86544423775356),
  @(44411, 0, 8, 70.86544423775356),
  @(44412, 1, 9, 79.72362476747276),
  @(44413, 4, 11, 97.43998582691115),
  @(44414, 1, 7, 62.00726370803437),
  @(44415, 4, 11, 97.43998582691115),
  @(44416, 6, 16, 141.7308884755071),
  @(44417, 2, 18, 159.4472495349455)
)

$startRow = 329

# Copy formatting from the last existing data row (328) down across the new rows first
$srcFormatRange = $ws.Range("A328:D328")
$endRow = $startRow + $data.Count - 1
$destFormatRange = $ws.Range("A$startRow`:D$endRow")
$srcFormatRange.Copy()
$destFormatRange.PasteSpecial(-4122)  # xlPasteFormats

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$wb.Save()
